$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column C to fit the new, longer "Preprocessing" note ---
$ws.Columns.Item(3).ColumnWidth = 20.6666667

# --- Add the new experiment row (row 8) ---
# Write the cells that introduce new shared strings first, in the same
# order they end up appended to sharedStrings.xml (G8, then C8, then E8).
$ws.Range("G8").Value = "accuracy    f1_macro    precision-neg    recall-neg" + [char]10 + "----------  ----------  ---------------  ------------" + [char]10 + "88.05%      88.05%      87.62%           88.59%"
$ws.Range("C8").Value = "lowercase" + [char]10 + "truncate 382/128 split"
$ws.Range("E8").Value = "seed = 1234" + [char]10 + "test_size = 0.2" + [char]10 + "MAX_LEN = 512" + [char]10 + "START_LEN = 382" + [char]10 + "END_LEN = 128" + [char]10 + "batch_size = 16" + [char]10 + "epochs = 10" + [char]10 + "use_gpu_test = True" + [char]10 + "lr = 1e-5"

$ws.Range("C8").WrapText = $true
$ws.Range("E8").WrapText = $true
$ws.Range("G8").WrapText = $true

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "bert-base-uncased"
$ws.Range("D8").Value = "NIL"
$ws.Range("F8").Value = 0.9

# Row height auto-grows with the wrapped, multi-line text (9 lines in E8)
$ws.Rows.Item(8).RowHeight = 129.6

# --- View changes: freeze header row, zoom to 85%, select F8 ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.Zoom = 85
$ws.Range("F8").Select() | Out-Null
